# Applies the edits described by the commit:
# - drop the KitchenCards sheet
# - rename MerchantCards -> Merchant, MixtureCards -> Mixture
# - update Merchant's "card type" (G) column header/values and widen the column
# - move the active-cell selections on the two remaining sheets

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$merchant = $wb.Worksheets.Item("MerchantCards")
$mixture  = $wb.Worksheets.Item("MixtureCards")
$kitchen  = $wb.Worksheets.Item("KitchenCards")

# Remove the KitchenCards sheet entirely.
$kitchen.Delete()

# Rename the remaining sheets.
$merchant.Name = "Merchant"
$mixture.Name = "Mixture"

# The merchant-hand marker cells ("*") become a descriptive placeholder string,
# centered like the rest of the header-aligned cells in that column.
$merchant.Range("G2").Value = "list<Mechant> inHand"
$merchant.Range("G2").HorizontalAlignment = -4108   # xlCenter
$merchant.Range("G2").VerticalAlignment = -4108     # xlCenter

$merchant.Range("G3").Value = "list<Mechant> inHand"
$merchant.Range("G3").HorizontalAlignment = -4108   # xlCenter
$merchant.Range("G3").VerticalAlignment = -4108     # xlCenter

# Column G needs to be widened to fit the new text.
$merchant.Columns.Item(7).ColumnWidth = 18.5703125

# Update the saved selections.
$merchant.Range("G7").Select()
$mixture.Range("A2").Select()
$mixture.Application.ActiveWindow.ScrollRow = 2

$merchant.Application.ActiveWindow.ScrollRow = 2
